# Apply "Added paq8l results on images" commit:
#  - new column Q ("PAQ8L -5") with per-image PAQ8L results
#  - column U ("Best") recomputed as MIN(C:T) now that Q participates
#  - misc view-state bookkeeping (best effort)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New column header -------------------------------------------------
$ws.Cells.Item(1, 17).Value = "PAQ8L -5"

# --- 2. New column Q data (PAQ8L -5 results), one value per existing data row
$q = @{
    2  = 104844
    3  = 111421
    5  = 74341
    6  = 81598
    8  = 133231
    9  = 137030
    11 = 165481
    12 = 170451
    14 = 101749
    15 = 110141
    17 = 419751
    18 = 434828
    20 = 42593
    21 = 44774
    23 = 150892
    24 = 157128
    26 = 161751
    27 = 167761
    29 = 418559
    30 = 369978
    32 = 168377
    33 = 175107
    35 = 93893
    36 = 99407
    38 = 90057
    39 = 97288
    41 = 161060
    42 = 165515
    44 = 186724
    45 = 191005
    47 = 151350
    48 = 162092
    50 = 121024
    51 = 125600
}

foreach ($r in $q.Keys) {
    $ws.Cells.Item($r, 17).Value2 = $q[$r]
}

# --- 3. Column U ("Best") becomes a real MIN() formula over C:T, so it
#        automatically folds in the new PAQ8L column and recalculates.
foreach ($r in $q.Keys) {
    $ws.Cells.Item($r, 21).Formula = "=MIN(C" + $r + ":T" + $r + ")"
}

# --- 4. Extend the duplicate-values highlight rule to cover the new U formulas
#        (mirrors the same rule/colors already used by the other highlight
#        rules on this sheet; added as its own rule since this runtime
#        collapses multi-area ModifyAppliesToRange calls to one area)
$newRange = $ws.Range("U3:U51")
$dupRule = $newRange.FormatConditions.AddUniqueValues()
$dupRule.DupeUnique = 1
$dupRule.Font.Color = 393372
$dupRule.Interior.Color = 13551615

# --- 5. Misc view-state bookkeeping (best effort) --------------------------
$ws.Range("Q1:Q1048576").Select() | Out-Null
